$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 - new "Investor 1" test offer row
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = 10
$ws.Range("A6").Style = "Normal 2"

$ws.Range("B6").Value = "Investor"
$ws.Range("B6").Style = "Normal 2"

$ws.Range("C6").Value = "Investor 1"

$ws.Range("D6").Value = "emp1@investor1.com"
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:emp1@investor1.com")
$ws.Range("D6").Style = "Hyperlink"

$ws.Range("E6").Value = "Emp1"
$ws.Range("F6").Value = "Investor1"

$ws.Range("G6").Value = "4415 Daniel Expressway, Doretheaside, IL 76585"
$ws.Range("G6").Style = "Normal 2"

$ws.Range("H6").Value = "UOUATXYOWU"
$ws.Range("H6").Style = "Normal 2"

$ws.Range("I6").Value = 209989880
$ws.Range("I6").Style = "Normal 2"

$ws.Range("J6").Value = "AX1123MM"
$ws.Range("J6").Style = "Normal 2"

$ws.Range("K6").Value = 5555

$ws.Range("L6").Value = "Bangalore"
$ws.Range("L6").Style = "Normal 2"

# ---------------------------------------------------------------------------
# Row 7 - new "Investor 2" test offer row
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 10
$ws.Range("A7").Style = "Normal 2"

$ws.Range("B7").Value = "Investor"
$ws.Range("B7").Style = "Normal 2"

$ws.Range("C7").Value = "Investor 2"

$ws.Range("D7").Value = "emp1@investor2.com"
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:emp1@investor2.com")
$ws.Range("D7").Style = "Hyperlink"

$ws.Range("E7").Value = "Emp1"
$ws.Range("F7").Value = "Investor2"

$ws.Range("G7").Value = "4416 Daniel Expressway, Doretheaside, IL 76585"
$ws.Range("G7").Style = "Normal 2"

$ws.Range("H7").Value = "AALLOPYT"
$ws.Range("H7").Style = "Normal 2"

$ws.Range("I7").Value = 8273487234
$ws.Range("I7").Style = "Normal 2"

$ws.Range("J7").Value = "SFDGSDFG"
$ws.Range("J7").Style = "Normal 2"

$ws.Range("K7").Value = 6666

$ws.Range("L7").Value = "Bangalore"
$ws.Range("L7").Style = "Normal 2"

# ---------------------------------------------------------------------------
# Move the active selection the way the author left it
# ---------------------------------------------------------------------------
$ws.Range("C8").Select()
